$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (column G) values, regenerated from std/mean calc (s_vals)
$kValues = @{
    2  = 1
    3  = 5
    4  = 5
    5  = 4
    6  = 4
    7  = 0
    8  = 0
    9  = 2
    10 = 1
    11 = 2
    12 = 0
    13 = 0
    14 = 3
    15 = 1
    16 = 0
    17 = 2
    18 = 3
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 1
    25 = 1
    26 = 1
    27 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
